$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

# Force a genuine content mutation first so the writer actually rebuilds
# the paragraph's runs (setting the identical concatenated text directly
# would be a no-op against the cached "Below section-level" read value
# and the three original runs would be left untouched).
$textRange.Text = "placeholder-to-force-run-rebuild"
$textRange.Text = "Below section-level"
